$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as forced text (preserves "t=str" style cell content
# without leaving a residual custom number-format style behind).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "274.89"
Set-TextValue "E2" "2.34%"
Set-TextValue "D3" "26.83"
Set-TextValue "E3" "0.73%"
Set-TextValue "D4" "4.926"
Set-TextValue "E4" "5.03%"
Set-TextValue "D5" "0.06362"
Set-TextValue "E5" "4.59%"
Set-TextValue "D6" "6.967"
Set-TextValue "E6" "3.37%"
Set-TextValue "D7" "3.362"
Set-TextValue "E7" "6.20%"
Set-TextValue "D8" "1.434"
Set-TextValue "E8" "61.06%"
Set-TextValue "D9" "0.8890"
Set-TextValue "E9" "4.75%"
Set-TextValue "D10" "0.1472"
Set-TextValue "E10" "3.86%"
Set-TextValue "D11" "0.05189"
Set-TextValue "E11" "6.26%"
Set-TextValue "D12" "0.07428"
Set-TextValue "E12" "4.78%"
Set-TextValue "D13" "0.03154"
Set-TextValue "E13" "-1.29%"
Set-TextValue "D14" "0.09067"
Set-TextValue "E14" "0.62%"
Set-TextValue "D15" "0.001565"
Set-TextValue "E15" "2.37%"
Set-TextValue "D16" "0.0006321"
Set-TextValue "E16" "3.94%"
Set-TextValue "D17" "0.006058"
Set-TextValue "E17" "-0.24%"
Set-TextValue "D18" "3.487"
Set-TextValue "E18" "0.89%"
Set-TextValue "D19" "2.281"
Set-TextValue "E19" "1.75%"
Set-TextValue "E20" "2.25%"
Set-TextValue "D21" "0.1336"
Set-TextValue "E21" "2.83%"
Set-TextValue "D22" "3.937"
Set-TextValue "E22" "2.38%"
Set-TextValue "D23" "0.04353"
Set-TextValue "E23" "3.03%"
Set-TextValue "D24" "0.001184"
Set-TextValue "E24" "0.33%"
Set-TextValue "D25" "0.003664"
Set-TextValue "E25" "-11.52%"
Set-TextValue "D26" "0.0001205"
Set-TextValue "E26" "0.55%"
Set-TextValue "D27" "0.0001945"
Set-TextValue "E27" "15.86%"
Set-TextValue "D40" "0.04041"
Set-TextValue "E40" "2.50%"
Set-TextValue "D41" "0.006655"
Set-TextValue "E41" "58.02%"
Set-TextValue "D42" "0.1170"
Set-TextValue "E42" "4.92%"
Set-TextValue "D43" "0.002370"
Set-TextValue "E43" "18.06%"
Set-TextValue "D44" "0.01223"
Set-TextValue "E44" "-2.60%"
Set-TextValue "D45" "0.00005248"
Set-TextValue "E45" "2.33%"
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D46" "2.355"
Set-TextValue "E46" "754.42%"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue "D47" "0.02129"
Set-TextValue "E47" "-13.03%"
$ws.Range("B48").Value = "SpecialPowerGold"
$ws.Range("C48").Value = "https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
Set-TextValue "D48" "0.0001999"
Set-TextValue "E48" "0.06%"
$ws.Range("B49").Value = "DigiFinexToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/rY6dWXQL4+digifinextoken-dft"
Set-TextValue "D49" "--"
Set-TextValue "E49" "--%"
$ws.Range("B50").Value = "Bitcoin2.0"
$ws.Range("C50").Value = "https://coinranking.com/coin/tSq1ehUma+bitcoin20-xbtc2"
Set-TextValue "D50" "--"
Set-TextValue "E50" "--%"
$ws.Range("B51").Value = "CoinField"
$ws.Range("C51").Value = "https://coinranking.com/coin/h4GpuIkN_+coinfield-cfc"
